# Updates TPM-derived NATMI LR-pair metrics (Col18a1-Ptprs) with recomputed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 40.34291466666667
$ws.Range("H2").Value = 121.028744
$ws.Range("I2").Value = 0.3404392602027053
$ws.Range("J2").Value = 0.3404392602027053
$ws.Range("M2").Value = 5.273684
$ws.Range("N2").Value = 15.821052
$ws.Range("O2").Value = 0.0510821201937383
$ws.Range("P2").Value = 0.0510821201937383
$ws.Range("Q2").Value = 212.7557835909653
$ws.Range("R2").Value = 1914.802052318688
$ws.Range("S2").Value = 0.01739035920834194
$ws.Range("T2").Value = 0.01739035920834194
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 40.34291466666667
$ws.Range("H3").Value = 121.028744
$ws.Range("I3").Value = 0.3404392602027053
$ws.Range("J3").Value = 0.3404392602027053
$ws.Range("O3").Value = 0.5598845502029881
$ws.Range("P3").Value = 0.5598845502029881
$ws.Range("Q3").Value = 2331.905483702174
$ws.Range("R3").Value = 20987.14935331957
$ws.Range("S3").Value = 0.1906066820700297
$ws.Range("T3").Value = 0.1906066820700297
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 40.34291466666667
$ws.Range("H4").Value = 121.028744
$ws.Range("I4").Value = 0.3404392602027053
$ws.Range("J4").Value = 0.3404392602027053
$ws.Range("M4").Value = 32.95839133333334
$ws.Range("N4").Value = 98.87517400000002
$ws.Range("O4").Value = 0.3192425840231603
$ws.Range("P4").Value = 0.3192425840231604
$ws.Range("Q4").Value = 1329.637569111273
$ws.Range("R4").Value = 11966.73812200146
$ws.Range("S4").Value = 0.1086827091300447
$ws.Range("T4").Value = 0.1086827091300447
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 40.34291466666667
$ws.Range("H5").Value = 121.028744
$ws.Range("I5").Value = 0.3404392602027053
$ws.Range("J5").Value = 0.3404392602027053
$ws.Range("M5").Value = 7.205150000000001
$ws.Range("N5").Value = 21.61545
$ws.Range("O5").Value = 0.06979074558011317
$ws.Range("P5").Value = 0.06979074558011318
$ws.Range("Q5").Value = 290.6767516105334
$ws.Range("R5").Value = 2616.090764494801
$ws.Range("S5").Value = 0.02375950979428895
$ws.Range("T5").Value = 0.02375950979428896
$ws.Range("I6").Value = 0.3497297648481489
$ws.Range("J6").Value = 0.3497297648481489
$ws.Range("M6").Value = 5.273684
$ws.Range("N6").Value = 15.821052
$ws.Range("O6").Value = 0.0510821201937383
$ws.Range("P6").Value = 0.0510821201937383
$ws.Range("Q6").Value = 218.561837201292
$ws.Range("R6").Value = 1967.056534811628
$ws.Range("S6").Value = 0.01786493788330097
$ws.Range("T6").Value = 0.01786493788330097
$ws.Range("I7").Value = 0.3497297648481489
$ws.Range("J7").Value = 0.3497297648481489
$ws.Range("O7").Value = 0.5598845502029881
$ws.Range("P7").Value = 0.5598845502029881
$ws.Range("S7").Value = 0.1958082920846026
$ws.Range("T7").Value = 0.1958082920846026
$ws.Range("I8").Value = 0.3497297648481489
$ws.Range("J8").Value = 0.3497297648481489
$ws.Range("M8").Value = 32.95839133333334
$ws.Range("N8").Value = 98.87517400000002
$ws.Range("O8").Value = 0.3192425840231603
$ws.Range("P8").Value = 0.3192425840231604
$ws.Range("Q8").Value = 1365.923055119054
$ws.Range("R8").Value = 12293.30749607149
$ws.Range("S8").Value = 0.1116486338399353
$ws.Range("T8").Value = 0.1116486338399353
$ws.Range("I9").Value = 0.3497297648481489
$ws.Range("J9").Value = 0.3497297648481489
$ws.Range("M9").Value = 7.205150000000001
$ws.Range("N9").Value = 21.61545
$ws.Range("O9").Value = 0.06979074558011317
$ws.Range("P9").Value = 0.06979074558011318
$ws.Range("Q9").Value = 298.60924949445
$ws.Range("R9").Value = 2687.48324545005
$ws.Range("S9").Value = 0.02440790104030996
$ws.Range("T9").Value = 0.02440790104030997
$ws.Range("G10").Value = 36.642055
$ws.Range("H10").Value = 109.926165
$ws.Range("I10").Value = 0.3092090445020276
$ws.Range("J10").Value = 0.3092090445020277
$ws.Range("M10").Value = 5.273684
$ws.Range("N10").Value = 15.821052
$ws.Range("O10").Value = 0.0510821201937383
$ws.Range("P10").Value = 0.0510821201937383
$ws.Range("Q10").Value = 193.23861918062
$ws.Range("R10").Value = 1739.14757262558
$ws.Range("S10").Value = 0.01579505357624355
$ws.Range("T10").Value = 0.01579505357624355
$ws.Range("G11").Value = 36.642055
$ws.Range("H11").Value = 109.926165
$ws.Range("I11").Value = 0.3092090445020276
$ws.Range("J11").Value = 0.3092090445020277
$ws.Range("O11").Value = 0.5598845502029881
$ws.Range("P11").Value = 0.5598845502029881
$ws.Range("Q11").Value = 2117.987996023903
$ws.Range("R11").Value = 19061.89196421513
$ws.Range("S11").Value = 0.1731213667997135
$ws.Range("T11").Value = 0.1731213667997135
$ws.Range("G12").Value = 36.642055
$ws.Range("H12").Value = 109.926165
$ws.Range("I12").Value = 0.3092090445020276
$ws.Range("J12").Value = 0.3092090445020277
$ws.Range("M12").Value = 32.95839133333334
$ws.Range("N12").Value = 98.87517400000002
$ws.Range("O12").Value = 0.3192425840231603
$ws.Range("P12").Value = 0.3192425840231604
$ws.Range("Q12").Value = 1207.663187947524
$ws.Range("R12").Value = 10868.96869152771
$ws.Range("S12").Value = 0.09871269437015967
$ws.Range("T12").Value = 0.0987126943701597
$ws.Range("G13").Value = 36.642055
$ws.Range("H13").Value = 109.926165
$ws.Range("I13").Value = 0.3092090445020276
$ws.Range("J13").Value = 0.3092090445020277
$ws.Range("M13").Value = 7.205150000000001
$ws.Range("N13").Value = 21.61545
$ws.Range("O13").Value = 0.06979074558011317
$ws.Range("P13").Value = 0.06979074558011318
$ws.Range("Q13").Value = 264.01150258325
$ws.Range("R13").Value = 2376.10352324925
$ws.Range("S13").Value = 0.0215799297559109
$ws.Range("T13").Value = 0.02157992975591091
$ws.Range("E14").Value = 3.0
$ws.Range("F14").Value = 1.0
$ws.Range("G14").Value = 0.07370033333333333
$ws.Range("H14").Value = 0.221101
$ws.Range("I14").Value = 0.0006219304471182344
$ws.Range("J14").Value = 0.0006219304471182345
$ws.Range("M14").Value = 5.273684
$ws.Range("N14").Value = 15.821052
$ws.Range("O14").Value = 0.0510821201937383
$ws.Range("P14").Value = 0.0510821201937383
$ws.Range("Q14").Value = 0.3886722686946666
$ws.Range("R14").Value = 3.498050418252
$ws.Range("S14").Value = 0.00003176952585183905
$ws.Range("T14").Value = 0.00003176952585183906
$ws.Range("E15").Value = 3.0
$ws.Range("F15").Value = 1.0
$ws.Range("G15").Value = 0.07370033333333333
$ws.Range("H15").Value = 0.221101
$ws.Range("I15").Value = 0.0006219304471182344
$ws.Range("J15").Value = 0.0006219304471182345
$ws.Range("O15").Value = 0.5598845502029881
$ws.Range("P15").Value = 0.5598845502029881
$ws.Range("Q15").Value = 4.260034577835778
$ws.Range("R15").Value = 38.340311200522
$ws.Range("S15").Value = 0.0003482092486423359
$ws.Range("T15").Value = 0.000348209248642336
$ws.Range("E16").Value = 3.0
$ws.Range("F16").Value = 1.0
$ws.Range("G16").Value = 0.07370033333333333
$ws.Range("H16").Value = 0.221101
$ws.Range("I16").Value = 0.0006219304471182344
$ws.Range("J16").Value = 0.0006219304471182345
$ws.Range("M16").Value = 32.95839133333334
$ws.Range("N16").Value = 98.87517400000002
$ws.Range("O16").Value = 0.3192425840231603
$ws.Range("P16").Value = 0.3192425840231604
$ws.Range("Q16").Value = 2.429044427397111
$ws.Range("R16").Value = 21.861399846574
$ws.Range("S16").Value = 0.0001985466830207046
$ws.Range("T16").Value = 0.0001985466830207047
$ws.Range("E17").Value = 3.0
$ws.Range("F17").Value = 1.0
$ws.Range("G17").Value = 0.07370033333333333
$ws.Range("H17").Value = 0.221101
$ws.Range("I17").Value = 0.0006219304471182344
$ws.Range("J17").Value = 0.0006219304471182345
$ws.Range("M17").Value = 7.205150000000001
$ws.Range("N17").Value = 21.61545
$ws.Range("O17").Value = 0.06979074558011317
$ws.Range("P17").Value = 0.06979074558011318
$ws.Range("Q17").Value = 0.5310219567166666
$ws.Range("R17").Value = 4.779197610450001
$ws.Range("S17").Value = 0.00004340498960335472
$ws.Range("T17").Value = 0.00004340498960335474

Write-Host "Applied 190 cell updates"
